# Added quiz seed data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C (questions) widened to fit the new, longer question text ---
$ws.Columns.Item(3).ColumnWidth = 27.666666666666668

# --- Quiz 1 (rows 1-5) ---
$ws.Range("C1").Value = "What planet is this?"
$ws.Range("C2").Value = "What is the name of our galaxy?"
$ws.Range("C3").Value = "In Geocentrism, which planet is the center of the galaxy?"
$ws.Range("D3").Value = "c"
$ws.Range("E3").ClearContents()
$ws.Range("C4").Value = "How many constellations do we recognize today?"
$ws.Range("C5").Value = "What constellation is this?"
$ws.Range("E5").Value = "Question5.png"

# --- Quiz 2 (rows 6-10) ---
$ws.Range("C6").Value = "Why do we have seasons?"
$ws.Range("C7").Value = "What angle is the Earth's axis tipped at?"
$ws.Range("C9").Value = "What moon phase is this?"
$ws.Range("D9").Value = "b"
$ws.Range("C10").Value = "What moon phase is this?"
$ws.Range("D10").Value = "c"
$ws.Range("E9").Value = "Question8.png"
$ws.Range("E9").Font.Color = 0
$ws.Range("E10").Value = "Question9.png"
$ws.Range("E10").Font.Color = 0
$ws.Range("C8").Value = "What is the moving line that separates the light and dark sides of the moon?"

# --- Quiz 3 (rows 11-15) ---
$ws.Range("C12").Value = "What is the penumbra?"
$ws.Range("C11").Value = "What is the umbra?"
$ws.Range("C13").Value = "What is the corona?"
$ws.Range("D13").Value = "d"
$ws.Range("C15").Value = "What is this effect called?"
$ws.Range("E15").Value = "Question15.png"
$ws.Range("C14").Value = "What type of eclipse occurs when the moon cannot cover the entire sun?"
$ws.Range("D14").Value = "b"

# --- Selection / page setup to match the edited author state ---
$ws.Range("F15").Select()
$ws.PageSetup.Orientation = 1
